$d = $word.ActiveDocument

# Locate the "Postman (Testen und Debuggen der API)" paragraph and add a
# new paragraph right after it for the new "OpenAPI Generator" tool entry,
# reproducing the same multi-run split the author typed it with.
$target = $null
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Postman*") {
        $target = $cand
    }
}

$r = $target.Range
$r.InsertParagraphAfter()
$newIndex = $target.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)
$ir = $newPara.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r><w:t>OpenAPI Generator</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> (Generieren </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">von </w:t></w:r>' +
  '<w:r><w:t>leeren</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> API</w:t></w:r>' +
  '<w:r><w:t>s</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> und Client</w:t></w:r>' +
  '<w:r><w:t>s</w:t></w:r>' +
  '<w:r><w:t>)</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$ir.InsertXML($xml) | Out-Null
